$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 5787, 45972.95833333334),
    @(3, 5752, 45972.96875),
    @(4, 5664, 45972.97916666666),
    @(5, 5591, 45972.98958333334),
    @(6, 5619, 45973),
    @(7, 5572, 45973.01041666666),
    @(8, 5550, 45973.02083333334),
    @(9, 5492, 45973.03125),
    @(10, 5459, 45973.04166666666),
    @(11, 5455, 45973.05208333334),
    @(12, 5375, 45973.0625),
    @(13, 5373, 45973.07291666666),
    @(14, 5432, 45973.08333333334),
    @(15, 5408, 45973.09375),
    @(16, 5381, 45973.10416666666),
    @(17, 5408, 45973.11458333334),
    @(18, 5367, 45973.125),
    @(19, 5431, 45973.13541666666),
    @(20, 5482, 45973.14583333334),
    @(21, 5558, 45973.15625),
    @(22, 5543, 45973.16666666666),
    @(23, 5681, 45973.17708333334),
    @(24, 5810, 45973.1875),
    @(25, 5717, 45973.19791666666),
    @(26, 6030, 45973.20833333334),
    @(27, 6233, 45973.21875),
    @(28, 6395, 45973.22916666666),
    @(29, 6572, 45973.23958333334),
    @(30, 6826, 45973.25),
    @(31, 7003, 45973.26041666666),
    @(32, 7124, 45973.27083333334),
    @(33, 7167, 45973.28125),
    @(34, 7333, 45973.29166666666),
    @(35, 7456, 45973.30208333334),
    @(36, 7542, 45973.3125),
    @(37, 7496, 45973.32291666666),
    @(38, 7486, 45973.33333333334),
    @(39, 7418, 45973.34375),
    @(40, 7406, 45973.35416666666),
    @(41, 7294, 45973.36458333334),
    @(42, 7113, 45973.375),
    @(43, 7157, 45973.38541666666),
    @(44, 7097, 45973.39583333334),
)

foreach ($item in $data) {
    $r = $item[0]
    $a = $item[1]
    $b = $item[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}

$ws.Range("B36:B44").NumberFormat = $ws.Range("B2").NumberFormat
